# Auto-generated Excel COM-interop script applying the Faerie_Profits market-price refresh.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR leve-profit tables to match the scheduled market-data run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 999
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H62").Value = 7393.696
$ws.Range("I62").Value = 7411.1577
$ws.Range("J62").Value = 7310.75
$ws.Range("K62").Value = 7411.1577
$ws.Range("L62").Value = 7310.75
$ws.Range("M62").Value = -6787.1577
$ws.Range("N62").Value = -8558.75

$ws.Range("H65").Value = 7393.696
$ws.Range("I65").Value = 7411.1577
$ws.Range("J65").Value = 7310.75
$ws.Range("K65").Value = 37055.7885
$ws.Range("L65").Value = 36553.75
$ws.Range("M65").Value = -33935.7885
$ws.Range("N65").Value = -42793.75

$ws.Range("H113").Value = 3666.3333
$ws.Range("I113").Value = 3249.5
$ws.Range("J113").Value = 3999.8
$ws.Range("K113").Value = 3249.5
$ws.Range("L113").Value = 3999.8
$ws.Range("M113").Value = 4.5
$ws.Range("N113").Value = -10507.8

$ws.Range("H116").Value = 2240
$ws.Range("I116").Value = 2240
$ws.Range("K116").Value = 2240
$ws.Range("M116").Value = 1202

$ws.Range("H132").Value = 7818.826
$ws.Range("I132").Value = 2674.5122
$ws.Range("K132").Value = 8023.5366
$ws.Range("M132").Value = -5493.5366

$ws.Range("H137").Value = 2764.3333
$ws.Range("I137").Value = 2864.8572
$ws.Range("J137").Value = 2412.5
$ws.Range("K137").Value = 8594.571599999999
$ws.Range("L137").Value = 7237.5
$ws.Range("M137").Value = -6044.571599999999
$ws.Range("N137").Value = -12337.5

$ws.Range("H138").Value = 173507.27
$ws.Range("I138").Value = 53602.633
$ws.Range("J138").Value = 224133.67
$ws.Range("K138").Value = 160807.899
$ws.Range("L138").Value = 672401.01
$ws.Range("M138").Value = -155667.899
$ws.Range("N138").Value = -682681.01

$ws.Range("H141").Value = 3218.8635
$ws.Range("I141").Value = 2313.0588
$ws.Range("K141").Value = 6939.176399999999
$ws.Range("M141").Value = -1759.176399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6564.71
$ws.Range("I32").Value = 6572.672
$ws.Range("K32").Value = 6572.672
$ws.Range("M32").Value = -6285.672

$ws.Range("H45").Value = 2757.2222
$ws.Range("I45").Value = 2352.9565
$ws.Range("J45").Value = 5081.75
$ws.Range("K45").Value = 2352.9565
$ws.Range("L45").Value = 5081.75
$ws.Range("M45").Value = -1975.9565
$ws.Range("N45").Value = -5835.75

$ws.Range("H61").Value = 5045.032
$ws.Range("I61").Value = 5052.1113
$ws.Range("J61").Value = 4997.25
$ws.Range("K61").Value = 5052.1113
$ws.Range("L61").Value = 4997.25
$ws.Range("M61").Value = -4840.1113
$ws.Range("N61").Value = -5421.25

$ws.Range("H132").Value = 2433.6204
$ws.Range("I132").Value = 2399.5857
$ws.Range("J132").Value = 2698.3333
$ws.Range("K132").Value = 7198.757100000001
$ws.Range("L132").Value = 8094.999899999999
$ws.Range("M132").Value = -4668.757100000001
$ws.Range("N132").Value = -13154.9999

$ws.Range("H136").Value = 5045.032
$ws.Range("I136").Value = 5052.1113
$ws.Range("J136").Value = 4997.25
$ws.Range("K136").Value = 15156.3339
$ws.Range("L136").Value = 14991.75
$ws.Range("M136").Value = -12606.3339
$ws.Range("N136").Value = -20091.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3738.0334
$ws.Range("I105").Value = 3584.318
$ws.Range("J105").Value = 4160.75
$ws.Range("K105").Value = 3584.318
$ws.Range("L105").Value = 4160.75
$ws.Range("M105").Value = -1837.318
$ws.Range("N105").Value = -7654.75

$ws.Range("H134").Value = 5727.2383
$ws.Range("I134").Value = 2302.2
$ws.Range("J134").Value = 8840.909
$ws.Range("K134").Value = 6906.599999999999
$ws.Range("L134").Value = 26522.727
$ws.Range("M134").Value = -4371.599999999999
$ws.Range("N134").Value = -31592.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1825.8334
$ws.Range("I58").Value = 1844.4
$ws.Range("J58").Value = 1733
$ws.Range("K58").Value = 1844.4
$ws.Range("L58").Value = 1733
$ws.Range("M58").Value = -1641.4
$ws.Range("N58").Value = -2139

$ws.Range("H62").Value = 2595.75
$ws.Range("I62").Value = 1092.5
$ws.Range("K62").Value = 1092.5
$ws.Range("M62").Value = -468.5

$ws.Range("H65").Value = 2595.75
$ws.Range("I65").Value = 1092.5
$ws.Range("K65").Value = 5462.5
$ws.Range("M65").Value = -2342.5

$ws.Range("H132").Value = 1178494.2
$ws.Range("I132").Value = 1430451
$ws.Range("J132").Value = 2695.6667
$ws.Range("K132").Value = 4291353
$ws.Range("L132").Value = 8087.000100000001
$ws.Range("M132").Value = -4288823
$ws.Range("N132").Value = -13147.0001

$ws.Range("H134").Value = 2852.3333
$ws.Range("I134").Value = 1357.7916
$ws.Range("J134").Value = 5841.4165
$ws.Range("K134").Value = 4073.3748
$ws.Range("L134").Value = 17524.2495
$ws.Range("M134").Value = -1538.3748
$ws.Range("N134").Value = -22594.2495

$ws.Range("H136").Value = 1825.8334
$ws.Range("I136").Value = 1844.4
$ws.Range("J136").Value = 1733
$ws.Range("K136").Value = 5533.200000000001
$ws.Range("L136").Value = 5199
$ws.Range("M136").Value = -2983.200000000001
$ws.Range("N136").Value = -10299

$ws.Range("H141").Value = 324574.25
$ws.Range("J141").Value = 324574.25
$ws.Range("L141").Value = 324574.25
$ws.Range("N141").Value = -334934.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 656.3570999999999
$ws.Range("I50").Value = 346.25
$ws.Range("J50").Value = 1069.8334
$ws.Range("K50").Value = 1038.75
$ws.Range("L50").Value = 3209.5002
$ws.Range("M50").Value = -557.75
$ws.Range("N50").Value = -4171.5002

$ws.Range("H53").Value = 656.3570999999999
$ws.Range("I53").Value = 346.25
$ws.Range("J53").Value = 1069.8334
$ws.Range("K53").Value = 1038.75
$ws.Range("L53").Value = 3209.5002
$ws.Range("M53").Value = -557.75
$ws.Range("N53").Value = -4171.5002

$ws.Range("H59").Value = 4337.5
$ws.Range("I59").Value = 1175
$ws.Range("K59").Value = 3525
$ws.Range("M59").Value = -2985

$ws.Range("H113").Value = 925.875
$ws.Range("J113").Value = 933.3570999999999
$ws.Range("L113").Value = 2800.0713
$ws.Range("N113").Value = -7140.0713

$ws.Range("H129").Value = 63380.625
$ws.Range("J129").Value = 2991.6667
$ws.Range("L129").Value = 8975.000100000001
$ws.Range("N129").Value = -18975.0001

$ws.Range("H131").Value = 910456.2
$ws.Range("I131").Value = 1429616.9
$ws.Range("J131").Value = 1925
$ws.Range("K131").Value = 4288850.699999999
$ws.Range("L131").Value = 5775
$ws.Range("M131").Value = -4283810.699999999
$ws.Range("N131").Value = -15855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4998.25
$ws.Range("I70").Value = 4999.6665
$ws.Range("J70").Value = 4994
$ws.Range("K70").Value = 4999.6665
$ws.Range("L70").Value = 4994
$ws.Range("M70").Value = -4729.6665
$ws.Range("N70").Value = -5534

$ws.Range("H73").Value = 4998.25
$ws.Range("I73").Value = 4999.6665
$ws.Range("J73").Value = 4994
$ws.Range("K73").Value = 4999.6665
$ws.Range("L73").Value = 4994
$ws.Range("M73").Value = -4063.6665
$ws.Range("N73").Value = -6866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3816
$ws.Range("I7").Value = 3033.0557
$ws.Range("J7").Value = 6164.8335
$ws.Range("K7").Value = 3033.0557
$ws.Range("L7").Value = 6164.8335
$ws.Range("M7").Value = -2921.0557
$ws.Range("N7").Value = -6388.8335

$ws.Range("H16").Value = 2999.1667
$ws.Range("I16").Value = 2352.0908
$ws.Range("J16").Value = 4016
$ws.Range("K16").Value = 2352.0908
$ws.Range("L16").Value = 4016
$ws.Range("M16").Value = -2182.0908
$ws.Range("N16").Value = -4356

$ws.Range("H40").Value = 5596.7856
$ws.Range("I40").Value = 5575.1904
$ws.Range("J40").Value = 5661.5713
$ws.Range("K40").Value = 5575.1904
$ws.Range("L40").Value = 5661.5713
$ws.Range("M40").Value = -5439.1904
$ws.Range("N40").Value = -5933.5713

$ws.Range("H126").Value = 3816
$ws.Range("I126").Value = 3033.0557
$ws.Range("J126").Value = 6164.8335
$ws.Range("K126").Value = 9099.167099999999
$ws.Range("L126").Value = 18494.5005
$ws.Range("M126").Value = -6629.167099999999
$ws.Range("N126").Value = -23434.5005

$ws.Range("H132").Value = 2816.0876
$ws.Range("I132").Value = 2774.9773
$ws.Range("K132").Value = 8324.9319
$ws.Range("M132").Value = -5794.9319

$ws.Range("H136").Value = 4418.864
$ws.Range("I136").Value = 3562.8333
$ws.Range("K136").Value = 10688.4999
$ws.Range("M136").Value = -8138.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8179.933
$ws.Range("I136").Value = 9559
$ws.Range("K136").Value = 28677
$ws.Range("M136").Value = -26127
